# Revert "adding term 2.0 now utf-8" (commit 78ae09ef7158944b5ce8ba326bcf0fbc536c742d)
#
# The prior commit had:
#   - bumped Version 1.1.0 -> 2.0.0
#   - bumped Date to 2024-06-03T10:45:43+02:00
#   - replaced the Contact placeholder with the full contact string
#   - replaced the "descendent-of" value "C" with a UUID on the
#     "Include from FSIII" sheet
#   - added a new "Include from FSIII 2" worksheet carrying the old
#     ("C") values
#
# This script undoes all of that.

$wb = $excel.ActiveWorkbook

# Remove the worksheet that was added by the reverted commit.
$wsExtra = $wb.Worksheets.Item("Include from FSIII 2")
$wsExtra.Delete()

# Restore the Metadata sheet's Version / Date / Contact values.
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.1.0"
$wsMeta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$wsMeta.Range("B10").Value = "No display for ContactDetail"

# Restore the "descendent-of" value on the remaining include sheet.
$wsInclude = $wb.Worksheets.Item("Include from FSIII")
$wsInclude.Range("C2").Value = "C"

# Re-select the Metadata sheet so the saved workbook's active tab
# matches the original (deleting the last sheet shifts selection).
$wsMeta.Activate()
